$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenarioA")

$ws.Range("F8").Value = 18
$ws.Range("I8").Value = 8
$ws.Range("I9").Value = 8
$ws.Range("F10").Value = 15
$ws.Range("I10").Value = 8
$ws.Range("I11").Value = 8
$ws.Range("F12").Value = 15
$ws.Range("I12").Value = 8
$ws.Range("I13").Value = 8
$ws.Range("F14").Value = 7
$ws.Range("I14").Value = 8
$ws.Range("I15").Value = 8
$ws.Range("F16").Value = 79
$ws.Range("I16").Value = 8
$ws.Range("I17").Value = 8
$ws.Range("I18").Value = 8

$ws.Range("I9:I18").Select()
